$wb = $excel.ActiveWorkbook

# row 125 (ALC) - hunk 0
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 4627.2
$ws.Cells.Item(125, 9).Value = 2803.5
$ws.Cells.Item(125, 10).Value = 6711.4287
$ws.Cells.Item(125, 11).Value = 25231.5
$ws.Cells.Item(125, 12).Value = 60402.85830000001
$ws.Cells.Item(125, 13).Value = -22771.5
$ws.Cells.Item(125, 14).Value = -65322.85830000001

# row 133 (ALC) - hunk 1
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(133, 8).Value = 76107.164
$ws.Cells.Item(133, 10).Value = 76107.164
$ws.Cells.Item(133, 12).Value = 76107.164
$ws.Cells.Item(133, 14).Value = -86227.164

# row 134 (ALC) - hunk 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(134, 8).Value = 40404.918
$ws.Cells.Item(134, 10).Value = 40404.918
$ws.Cells.Item(134, 12).Value = 40404.918
$ws.Cells.Item(134, 14).Value = -50544.918

# row 136 (ALC) - hunk 3
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(136, 8).Value = 57579.5
$ws.Cells.Item(136, 10).Value = 57579.5
$ws.Cells.Item(136, 12).Value = 57579.5
$ws.Cells.Item(136, 14).Value = -67779.5

# row 138 (ALC) - hunk 4
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 1645.4688
$ws.Cells.Item(138, 9).Value = 1127.2778
$ws.Cells.Item(138, 10).Value = 2311.7144
$ws.Cells.Item(138, 11).Value = 3381.8334
$ws.Cells.Item(138, 12).Value = 6935.1432
$ws.Cells.Item(138, 13).Value = 1758.1666
$ws.Cells.Item(138, 14).Value = -17215.1432

# row 140 (ALC) - hunk 5
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(140, 8).Value = 63871.918
$ws.Cells.Item(140, 10).Value = 65068.547
$ws.Cells.Item(140, 12).Value = 65068.547
$ws.Cells.Item(140, 14).Value = -75428.54699999999

# row 5 (ARM) - hunk 6
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 259.5
$ws.Cells.Item(5, 9).Value = 259.5
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 259.5
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = -147.5
$ws.Cells.Item(5, 14).ClearContents()

# row 45 (ARM) - hunk 7
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 9618262
$ws.Cells.Item(45, 9).Value = 2482.2
$ws.Cells.Item(45, 11).Value = 2482.2
$ws.Cells.Item(45, 13).Value = -2105.2

# row 52 (ARM) - hunk 8
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(52, 8).Value = 55996.4
$ws.Cells.Item(52, 10).Value = 55996.4
$ws.Cells.Item(52, 12).Value = 55996.4
$ws.Cells.Item(52, 14).Value = -56632.4

# row 121 (ARM) - hunk 9
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(121, 8).Value = 48395.6
$ws.Cells.Item(121, 10).Value = 48395.6
$ws.Cells.Item(121, 12).Value = 48395.6
$ws.Cells.Item(121, 14).Value = -51889.6

# row 124 (ARM) - hunk 10
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(124, 8).Value = 51251.332
$ws.Cells.Item(124, 10).Value = 51251.332
$ws.Cells.Item(124, 12).Value = 51251.332
$ws.Cells.Item(124, 14).Value = -61071.332

# row 4 (BSM) - hunk 11
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 259.5
$ws.Cells.Item(4, 9).Value = 259.5
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 259.5
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = -144.5
$ws.Cells.Item(4, 14).ClearContents()

# row 51 (BSM) - hunk 12
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(51, 8).Value = 42985
$ws.Cells.Item(51, 10).Value = 42985
$ws.Cells.Item(51, 12).Value = 42985
$ws.Cells.Item(51, 14).Value = -43967

# row 52 (BSM) - hunk 13
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(52, 8).Value = 99990
$ws.Cells.Item(52, 10).Value = 99990
$ws.Cells.Item(52, 12).Value = 99990
$ws.Cells.Item(52, 14).Value = -100516

# row 94 (BSM) - hunk 14
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1841.2
$ws.Cells.Item(94, 9).Value = 2012.7778
$ws.Cells.Item(94, 11).Value = 2012.7778
$ws.Cells.Item(94, 13).Value = -1561.7778

# row 109 (BSM) - hunk 15
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(109, 8).Value = 78282.71000000001
$ws.Cells.Item(109, 10).Value = 78282.71000000001
$ws.Cells.Item(109, 12).Value = 78282.71000000001
$ws.Cells.Item(109, 14).Value = -81056.71000000001

# row 121 (BSM) - hunk 16
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(121, 8).Value = 99990
$ws.Cells.Item(121, 10).Value = 99990
$ws.Cells.Item(121, 12).Value = 99990
$ws.Cells.Item(121, 14).Value = -103484

# row 127 (BSM) - hunk 17
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(127, 8).Value = 59409
$ws.Cells.Item(127, 10).Value = 59409
$ws.Cells.Item(127, 12).Value = 59409
$ws.Cells.Item(127, 14).Value = -69329

# row 132 (BSM) - hunk 18
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(132, 8).Value = 32473.025
$ws.Cells.Item(132, 10).Value = 32473.025
$ws.Cells.Item(132, 12).Value = 32473.025
$ws.Cells.Item(132, 14).Value = -42593.025

# row 135 (BSM) - hunk 19
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(135, 8).Value = 72741.14
$ws.Cells.Item(135, 10).Value = 72741.14
$ws.Cells.Item(135, 12).Value = 72741.14
$ws.Cells.Item(135, 14).Value = -82881.14

# row 138 (BSM) - hunk 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(138, 8).Value = 68776.086
$ws.Cells.Item(138, 10).Value = 68776.086
$ws.Cells.Item(138, 12).Value = 68776.086
$ws.Cells.Item(138, 14).Value = -79056.086

# row 140 (BSM) - hunk 21
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(140, 8).Value = 43499
$ws.Cells.Item(140, 10).Value = 43499
$ws.Cells.Item(140, 12).Value = 43499
$ws.Cells.Item(140, 14).Value = -53859

# row 7 (CRP) - hunk 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 6905.8667
$ws.Cells.Item(7, 9).Value = 6855.2
$ws.Cells.Item(7, 10).Value = 6956.533
$ws.Cells.Item(7, 11).Value = 6855.2
$ws.Cells.Item(7, 12).Value = 6956.533
$ws.Cells.Item(7, 13).Value = -6742.2
$ws.Cells.Item(7, 14).Value = -7182.533

# row 22 (CRP) - hunk 23
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 497
$ws.Cells.Item(22, 9).Value = 321.5
$ws.Cells.Item(22, 10).Value = 1199
$ws.Cells.Item(22, 11).Value = 321.5
$ws.Cells.Item(22, 12).Value = 1199
$ws.Cells.Item(22, 13).Value = 28.5
$ws.Cells.Item(22, 14).Value = -1899

# row 93 (CRP) - hunk 24
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(93, 8).Value = 66333
$ws.Cells.Item(93, 9).Value = 66333
$ws.Cells.Item(93, 11).Value = 66333
$ws.Cells.Item(93, 13).Value = -64461

# row 138 (CRP) - hunk 25
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(138, 8).Value = 99996
$ws.Cells.Item(138, 10).Value = 99996
$ws.Cells.Item(138, 12).Value = 99996
$ws.Cells.Item(138, 14).Value = -110276

# row 5 (CUL) - hunk 26
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 952.7857
$ws.Cells.Item(5, 9).Value = 789.8333
$ws.Cells.Item(5, 11).Value = 2369.4999
$ws.Cells.Item(5, 13).Value = -2257.4999

# row 54 (CUL) - hunk 27
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(54, 8).Value = 3402
$ws.Cells.Item(54, 9).Value = 275
$ws.Cells.Item(54, 10).Value = 4444.3335
$ws.Cells.Item(54, 11).Value = 825
$ws.Cells.Item(54, 12).Value = 13333.0005
$ws.Cells.Item(54, 13).Value = -266
$ws.Cells.Item(54, 14).Value = -14451.0005

# row 55 (CUL) - hunk 28
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 1905.4445
$ws.Cells.Item(55, 9).Value = 714.1429000000001
$ws.Cells.Item(55, 10).Value = 6075
$ws.Cells.Item(55, 11).Value = 2142.4287
$ws.Cells.Item(55, 12).Value = 18225
$ws.Cells.Item(55, 13).Value = -1965.4287
$ws.Cells.Item(55, 14).Value = -18579

# row 135 (CUL) - hunk 29
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 952.7857
$ws.Cells.Item(135, 9).Value = 789.8333
$ws.Cells.Item(135, 11).Value = 7108.4997
$ws.Cells.Item(135, 13).Value = -4573.4997

# row 137 (CUL) - hunk 30
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value = 3479.0476
$ws.Cells.Item(137, 10).Value = 4635.5454
$ws.Cells.Item(137, 12).Value = 13906.6362
$ws.Cells.Item(137, 14).Value = -24106.6362

# row 109 (GSM) - hunk 31
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(109, 8).Value = 31999.143
$ws.Cells.Item(109, 10).Value = 31999.143
$ws.Cells.Item(109, 12).Value = 31999.143
$ws.Cells.Item(109, 14).Value = -34079.143

# row 113 (GSM) - hunk 32
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 3134249.5
$ws.Cells.Item(113, 9).Value = 224942.2
$ws.Cells.Item(113, 10).Value = 5558672.5
$ws.Cells.Item(113, 11).Value = 224942.2
$ws.Cells.Item(113, 12).Value = 5558672.5
$ws.Cells.Item(113, 13).Value = -222772.2
$ws.Cells.Item(113, 14).Value = -5563012.5

# row 119 (GSM) - hunk 33
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(119, 8).Value = 59997
$ws.Cells.Item(119, 10).Value = 59997
$ws.Cells.Item(119, 12).Value = 59997
$ws.Cells.Item(119, 14).Value = -69673

# row 123 (GSM) - hunk 34
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(123, 8).Value = 50959.2
$ws.Cells.Item(123, 10).Value = 50959.2
$ws.Cells.Item(123, 12).Value = 50959.2
$ws.Cells.Item(123, 14).Value = -55859.2

# row 135 (GSM) - hunk 35
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(135, 8).Value = 45436.25
$ws.Cells.Item(135, 10).Value = 45436.25
$ws.Cells.Item(135, 12).Value = 45436.25
$ws.Cells.Item(135, 14).Value = -55576.25

# row 140 (GSM) - hunk 36
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(140, 8).Value = 95206.664
$ws.Cells.Item(140, 10).Value = 95634.55
$ws.Cells.Item(140, 12).Value = 95634.55
$ws.Cells.Item(140, 14).Value = -105994.55

# row 96 (LTW) - hunk 37
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(96, 8).Value = 30000
$ws.Cells.Item(96, 10).Value = 30000
$ws.Cells.Item(96, 12).Value = 30000
$ws.Cells.Item(96, 14).Value = -35492

# row 122 (LTW) - hunk 38
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 30773044
$ws.Cells.Item(122, 9).Value = 37040416
$ws.Cells.Item(122, 10).Value = 16671452
$ws.Cells.Item(122, 11).Value = 111121248
$ws.Cells.Item(122, 12).Value = 50014356
$ws.Cells.Item(122, 13).Value = -111118798
$ws.Cells.Item(122, 14).Value = -50019256

# row 62 (WVR) - hunk 39
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 5650
$ws.Cells.Item(62, 9).Value = 5812.5
$ws.Cells.Item(62, 10).Value = 5000
$ws.Cells.Item(62, 11).Value = 5812.5
$ws.Cells.Item(62, 12).Value = 5000
$ws.Cells.Item(62, 13).Value = -5188.5
$ws.Cells.Item(62, 14).Value = -6248

# row 65 (WVR) - hunk 40
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(65, 8).Value = 5650
$ws.Cells.Item(65, 9).Value = 5812.5
$ws.Cells.Item(65, 10).Value = 5000
$ws.Cells.Item(65, 11).Value = 29062.5
$ws.Cells.Item(65, 12).Value = 25000
$ws.Cells.Item(65, 13).Value = -25942.5
$ws.Cells.Item(65, 14).Value = -31240

# row 94 (WVR) - hunk 41
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(94, 8).Value = 30312.5
$ws.Cells.Item(94, 9).Value = 19000
$ws.Cells.Item(94, 10).Value = 31928.572
$ws.Cells.Item(94, 11).Value = 19000
$ws.Cells.Item(94, 12).Value = 31928.572
$ws.Cells.Item(94, 13).Value = -18099
$ws.Cells.Item(94, 14).Value = -33730.572

# row 122 (WVR) - hunk 42
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1525.7693
$ws.Cells.Item(122, 9).Value = 1338.5333
$ws.Cells.Item(122, 10).Value = 1781.091
$ws.Cells.Item(122, 11).Value = 4015.5999
$ws.Cells.Item(122, 12).Value = 5343.272999999999
$ws.Cells.Item(122, 13).Value = -1565.5999
$ws.Cells.Item(122, 14).Value = -10243.273

# row 126 (WVR) - hunk 43
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 72779.69
$ws.Cells.Item(126, 9).Value = 103520.55
$ws.Cells.Item(126, 10).Value = 5149.8
$ws.Cells.Item(126, 11).Value = 310561.65
$ws.Cells.Item(126, 12).Value = 15449.4
$ws.Cells.Item(126, 13).Value = -308091.65
$ws.Cells.Item(126, 14).Value = -20389.4

# row 127 (WVR) - hunk 44
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(127, 8).Value = 0
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(127, 14).ClearContents()
